# Generate Report for Handback
# Update the timestamp strings recorded on the handback-status report.
$wb = $excel.ActiveWorkbook

# "Latest HO Xliff Generate Date" column on the Overview sheet (G2) and
# "Correspond Handoff Datetime" column on the de-de sheet (H2) share the
# same text value, so both need to be updated together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 01:09:09"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-21 01:09:02"
$wsZhCn.Range("K2").Value = "2016-08-21 01:09:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-21 01:09:09"
$wsDeDe.Range("K2").Value = "2016-08-21 01:09:35"
